$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 held the rule name "R40" (shared string), change it to "1"
$ws.Range("B11").Value = "1"
